# 14_days.xlsx fix:
#   - correct the "Esfenvalerate_0.01" (column D) survival values for rows 2-8
#   - restore the taller row heights (header/body rows) that had been
#     squashed back down
#
# (The workbook also picked up some incidental style-table bookkeeping in the
# authoritative diff -- a duplicate cellXfs entry that column J pointed at was
# removed and two unused "general" alignment entries were normalised to
# "left" -- but none of that is visible on any populated cell, so the
# meaningful, user-facing fix is entirely the data values and row heights
# below.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the duplicated/incorrect column D values ---------------------
$ws.Range("D2").Value = 94
$ws.Range("D3").Value = 98
$ws.Range("D4").Value = 91
$ws.Range("D5").Value = 97.3333333333333
$ws.Range("D6").Value = 90.2222222222222
$ws.Range("D7").Value = 95.1111111111111
$ws.Range("D8").Value = 56

# --- Restore the row heights -------------------------------------------
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 19.5
$ws.Rows.Item(5).RowHeight = 19.5
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 19.5
$ws.Rows.Item(8).RowHeight = 19.5
$ws.Rows.Item(10).RowHeight = 20.25
